$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Software Requirements paragraph: "20XX.X" -> "R2019b"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "MATLAB 20XX.X, but will",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "MATLAB R2019b, but will", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Introduction paragraph: remove the old hidden "_GoBack" bookmark
#    that currently sits right after " [link]".
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 3) Software Requirements paragraph: "[list]" -> actual package list
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "packages should be downloaded and installed: [list].  ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "packages should be downloaded and installed: Curve Fitting Toolbox, Statistics and Machine Learning Toolbox, and the Parallel Computing Toolbox (if parallelization is going to be used to speed up clustering).  ",
    2) | Out-Null

# ------------------------------------------------------------------
# 4) Re-insert the "_GoBack" bookmark right after "...clustering)" and
#    before the trailing ".  " that now ends that sentence.
# ------------------------------------------------------------------
$markRange = $d.Content
$markRange.Find.Execute(
    "speed up clustering)",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$markRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $markRange) | Out-Null
